# "include no rank decision in binary" - refresh the re_range scores (col H)
# for every worker and update age/prolificid/name (cols D/E/F) and race (col I)
# for the rows whose rank order shifted once the "no rank decision" cases were
# folded back into the binary-ranking computation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 7.360079283446961
$ws.Range("H3").Value = 6.358750456454161
$ws.Range("H4").Value = 6.017579664918089
$ws.Range("H5").Value = 5.309718579672998
$ws.Range("H6").Value = 5.244195657518464
$ws.Range("H7").Value = 4.045026469112039
$ws.Range("H8").Value = 1.242073243576292
$ws.Range("H9").Value = 1.014010395470444
$ws.Range("D10").Value = 30
$ws.Range("E10").Value = "60d5775a99b502eec8cf56b4"
$ws.Range("F10").Value = "Shadaisia"
$ws.Range("H10").Value = 0.4834459824271087
$ws.Range("D11").Value = 32
$ws.Range("E11").Value = "6036f9b3b1842f8b659b18c7"
$ws.Range("F11").Value = "Kellie"
$ws.Range("H11").Value = 0.1753989618967279
$ws.Range("I11").Value = "White"
$ws.Range("D12").Value = 33
$ws.Range("E12").Value = "60cb36ee9f58331a33cf5506"
$ws.Range("F12").Value = "Shaniek"
$ws.Range("H12").Value = 0.1348575153764161
$ws.Range("I12").Value = "Black or African American"
$ws.Range("D13").Value = 21
$ws.Range("E13").Value = "5c0e89c6c323400001e6c4a5"
$ws.Range("F13").Value = "Bri"
$ws.Range("H13").Value = 0.1140016948445168
$ws.Range("H14").Value = 13.17322371252606
$ws.Range("H15").Value = 8.25975769155853
$ws.Range("D16").Value = 27
$ws.Range("E16").Value = "5ff8ad350d084e10f500e48a"
$ws.Range("F16").Value = "Drew"
$ws.Range("H16").Value = 7.198709993617562
$ws.Range("D17").Value = 30
$ws.Range("E17").Value = "60c2341fe95d71ee52c043f0"
$ws.Range("F17").Value = "Matthew"
$ws.Range("H17").Value = 7.013420770724821
$ws.Range("H18").Value = 5.30118687809812
$ws.Range("H19").Value = 5.217927984380697
$ws.Range("H20").Value = 5.022934074744907
$ws.Range("D21").Value = 2
$ws.Range("E21").Value = "5e2522d6b734b47915f88275"
$ws.Range("F21").Value = "Corey"
$ws.Range("H21").Value = 4.346124827430741
$ws.Range("D22").Value = 33
$ws.Range("E22").Value = "60b322994d0b901954690036"
$ws.Range("F22").Value = "Brennan"
$ws.Range("H22").Value = 4.218944548501164
$ws.Range("H23").Value = 3.102299781206951
$ws.Range("D24").Value = 29
$ws.Range("E24").Value = "60b83826821417f8e484a207"
$ws.Range("F24").Value = "Eli"
$ws.Range("H24").Value = 2.368353172506046
$ws.Range("I24").Value = "White"
$ws.Range("D25").Value = 50
$ws.Range("E25").Value = "6097b95056caf5ebb2720002"
$ws.Range("F25").Value = "Damian"
$ws.Range("H25").Value = 2.301699192143767
$ws.Range("I25").Value = "Black or African American"
